# Re-applies the Aug 5 2023 GitHub Actions price refresh to the cryptos sheet:
#  - updates Price (D) / Volume(1h) (E) for every existing coin row
#  - inserts a new "BabyDogeCoin" row at position 48, which pushes every
#    row below it down by one (TheSandbox -> 49, XinFinNetwork -> 51,
#    Cronos drops off the bottom of the A1:E51 range)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Column D sometimes holds plain numeric-looking text (e.g. "0.9996").
    # Prefixing with an apostrophe (quote-prefix) is how Excel keeps typed
    # input as text instead of silently converting it to a Number, which
    # would drop significant trailing/leading zeros.
    $ws.Range($addr).Value = ('''' + $value)
}

$ws.Range('D2').Value = '29.018.82'
$ws.Range('E2').Value = '  -0.42%  '

$ws.Range('D3').Value = '1.830.98'
$ws.Range('E3').Value = '  -0.14%  '

Set-TextValue 'D4' '0.9996'
$ws.Range('E4').Value = '  +0.01%  '

Set-TextValue 'D5' '241.44'
$ws.Range('E5').Value = '  -0.12%  '

Set-TextValue 'D6' '0.6251'
$ws.Range('E6').Value = '  -5.51%  '

$ws.Range('E7').Value = '  +0.04%  '

Set-TextValue 'D8' '45.39'
$ws.Range('E8').Value = '  +2.22%  '

Set-TextValue 'D9' '0.07591'
$ws.Range('E9').Value = '  +2.47%  '

Set-TextValue 'D10' '0.2910'
$ws.Range('E10').Value = '  -0.89%  '

Set-TextValue 'D11' '22.69'
$ws.Range('E11').Value = '  -1.23%  '

Set-TextValue 'D12' '0.07747'
$ws.Range('E12').Value = '  +0.03%  '

$ws.Range('D13').Value = '1.830.38'
$ws.Range('E13').Value = '  -0.28%  '

Set-TextValue 'D14' '4.950'
$ws.Range('E14').Value = '  -0.98%  '

Set-TextValue 'D15' '0.6615'
$ws.Range('E15').Value = '  -1.21%  '

Set-TextValue 'D16' '82.40'
$ws.Range('E16').Value = '  -0.63%  '

Set-TextValue 'D17' '0.000009489'
$ws.Range('E17').Value = '  +10.15%  '

Set-TextValue 'D18' '5.959'
$ws.Range('E18').Value = '  -2.44%  '

$ws.Range('D19').Value = '29.003.57'
$ws.Range('E19').Value = '  -0.53%  '

Set-TextValue 'D20' '223.85'
$ws.Range('E20').Value = '  -1.33%  '

Set-TextValue 'D21' '12.31'
$ws.Range('E21').Value = '  -1.35%  '

$ws.Range('E22').Value = '  -0.01%  '

Set-TextValue 'D23' '7.191'
$ws.Range('E23').Value = '  +0.80%  '

$ws.Range('E24').Value = '  +0.04%  '

Set-TextValue 'D25' '159.84'
$ws.Range('E25').Value = '  +0.50%  '

Set-TextValue 'D26' '8.410'
$ws.Range('E26').Value = '  -2.17%  '

Set-TextValue 'D27' '0.1360'
$ws.Range('E27').Value = '  -3.33%  '

Set-TextValue 'D28' '17.79'
$ws.Range('E28').Value = '  -1.26%  '

Set-TextValue 'D29' '1.494'
$ws.Range('E29').Value = '  -1.39%  '

Set-TextValue 'D30' '4.053'
$ws.Range('E30').Value = '  -1.59%  '

Set-TextValue 'D31' '4.020'
$ws.Range('E31').Value = '  -0.86%  '

$ws.Range('E32').Value = '  +0.65%  '

Set-TextValue 'D33' '0.05183'
$ws.Range('E33').Value = '  -1.80%  '

Set-TextValue 'D34' '1.843'
$ws.Range('E34').Value = '  -1.46%  '

Set-TextValue 'D35' '0.7348'
$ws.Range('E35').Value = '  -0.50%  '

Set-TextValue 'D36' '1.146'
$ws.Range('E36').Value = '  -0.07%  '

$ws.Range('E37').Value = '  +2.04%  '

$ws.Range('D38').Value = '1.262.55'
$ws.Range('E38').Value = '  -2.95%  '

Set-TextValue 'D39' '2.757'
$ws.Range('E39').Value = '  +0.57%  '

Set-TextValue 'D40' '0.01783'
$ws.Range('E40').Value = '  -0.30%  '

Set-TextValue 'D41' '6.277'
$ws.Range('E41').Value = '  +0.03%  '

Set-TextValue 'D42' '0.8905'
$ws.Range('E42').Value = '  -2.82%  '

Set-TextValue 'D43' '1.000'
$ws.Range('E43').Value = '  +0.13%  '

Set-TextValue 'D44' '101.56'
$ws.Range('E44').Value = '  -0.82%  '

$ws.Range('D45').Value = '1.976.75'
$ws.Range('E45').Value = '  -0.63%  '

Set-TextValue 'D46' '64.43'
$ws.Range('E46').Value = '  +0.82%  '

$ws.Range('E47').Value = '  -0.45%  '

$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D48' '0.00000000120'
$ws.Range('E48').Value = '  +0.30%  '

$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D49' '0.3975'
$ws.Range('E49').Value = '  -0.74%  '

Set-TextValue 'D50' '8.825'
$ws.Range('E50').Value = '  -0.18%  '

$ws.Range('B51').Value = 'XinFinNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
Set-TextValue 'D51' '0.07149'
$ws.Range('E51').Value = '  -13.39%  '
